$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the color code cells D2, E2, D4, E4 from numeric 0 to the text "black"
$ws.Range("D2").Value = "black"
$ws.Range("E2").Value = "black"
$ws.Range("D4").Value = "black"
$ws.Range("E4").Value = "black"

# Update the active selection to D4
$ws.Range("D4").Select()
